# Add a "Save" column (H) to the s_vals sheet, matching the style of the
# existing header row (G1, e.g. "sum") and filling in the boolean-like
# 0/1 values for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, border, centered/top alignment)
# from the last existing header cell (G1) onto the new header cell (H1)
# so it reuses the same cell style rather than creating a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 0
